$d = $word.ActiveDocument

# 1. Remove the "Location: 631 Moston Lane, M40 5QD" paragraph entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Location: 631 Moston Lane, M40 5QD*") {
        $p.Range.Delete()
        break
    }
}
